$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4, 6, 7, 8, 9)
$values = @{
    2 = 0.4723336674725943
    3 = 0.5149173457127998
    4 = 0.7933767163891069
    6 = 5.263405291768424
    7 = 9.337773373270942
    8 = 25.85684824586824
    9 = 212.5690296848721
}

foreach ($r in $rows) {
    $ws.Range("B$r").Value = $values[$r]
    $ws.Range("C$r").Formula = "=AVERAGE(B$r`:B$r)"
    $ws.Range("D$r").Formula = "=STDEV(B$r`:B$r)"
}

$wb.Save()
